# GitHub Actions symbol-list refresh (Thu Dec 15 15:41:40 UTC 2022)
# Updates coin prices and re-ranks a few coins whose relative ordering
# changed, per the scraped coinranking.com snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: prices are stored as literal text in this sheet (t="inlineStr"),
# so every numeric-looking update is written with a leading apostrophe to
# force Excel to keep it as text instead of silently coercing to a Number.

# --- Price-only updates (rows whose coin/rank did not change) ---
$ws.Range("D2").Value  = "'264.37"
$ws.Range("D3").Value  = "'22.75"
$ws.Range("D4").Value  = "'6.206"
$ws.Range("D5").Value  = "'0.06088"
$ws.Range("D6").Value  = "'3.523"
$ws.Range("D7").Value  = "'6.722"
$ws.Range("D8").Value  = "'1.365"
$ws.Range("D9").Value  = "'0.8145"
$ws.Range("D10").Value = "'0.1589"
$ws.Range("D11").Value = "'0.08157"
$ws.Range("D12").Value = "'0.03364"
$ws.Range("D13").Value = "'0.03168"
$ws.Range("D14").Value = "'0.09257"
$ws.Range("D15").Value = "'3.932"
$ws.Range("D16").Value = "'0.001695"
$ws.Range("D17").Value = "'0.04844"
$ws.Range("D18").Value = "'0.0006238"
$ws.Range("D19").Value = "'0.006191"
$ws.Range("D20").Value = "'0.006107"
$ws.Range("D21").Value = "'0.001100"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("D23").Value = "'3.695"
$ws.Range("D24").Value = "'2.282"
$ws.Range("D25").Value = "'0.3385"
$ws.Range("D26").Value = "'0.1269"
$ws.Range("D27").Value = "'0.0002685"
$ws.Range("D40").Value = "'0.04635"

# --- Rows 41-43 swap rank order: KickToken/BKEXToken/CEJI -> BKEXToken/CEJI/KickToken ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1122"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003135"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003455"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Remaining price-only updates ---
$ws.Range("D44").Value = "'0.01053"
$ws.Range("D45").Value = "'0.00006092"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.7510"

# Row 48 (BOLO): price jump + "Best in 24h" tag appended to volume label
$ws.Range("D48").Value = "'0.1776"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"

# Row 49 (CryptobidCoin): price tweak + "Best in 24h" tag removed from volume label
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("E49").Value = "48CryptobidCoinCBC"

$ws.Range("D50").Value = "'0.01242"
